$wb = $excel.ActiveWorkbook

# The workbook currently has 4 sheets: 2021-Q2, 2021-Q3, 2021-Q4, 总计 (this
# last one holds the "total" roll-up table). We need to:
#   1. Insert a new "2022-Q1" sheet (same per-quarter layout as the other
#      quarter sheets) right before the "总计" sheet.
#   2. Add a new row to the "总计" roll-up table for 2022-Q1 (at the top,
#      pushing the existing rows down).
#
# The existing "总计" sheet object is reused/renamed to become "2022-Q1" (so
# it keeps its original sheetId/relationship slot) and a fresh sheet is
# created after it to become the new "总计" sheet - this mirrors how the
# workbook was actually restructured.

$totalOld = $wb.Worksheets.Item("总计")

# --- Grab the "总计" sheet's existing formatting before we touch anything,
# so the new "总计" sheet can reuse the same header / index-column styles.
$newTotal = $wb.Worksheets.Add($null, $totalOld)
$totalOld.Range("B1:D1").Copy($newTotal.Range("B1:D1"))
$totalOld.Range("A2").Copy($newTotal.Range("A2:A5"))

# Rename sheets into their final names.
$totalOld.Name = "2022-Q1"
$newTotal.Name = "总计"

$q1 = $totalOld

# Drop the leftover roll-up rows (the old "总计" sheet had 4 data rows; the
# new "2022-Q1" per-fund sheet only needs 1).
$q1.Rows("3:4").Delete()

# --- Build the "2022-Q1" per-fund sheet (same columns as the other quarter
# sheets: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名).
# Extend the existing header style (already on B1:D1 from the old "总计"
# sheet) across the rest of the header row.
$q1.Range("D1").Copy($q1.Range("E1:H1"))

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
# Fund code / ratio-like figures are stored as plain text in this workbook
# (e.g. "004223" keeps its leading zero). Use the leading-apostrophe form to
# force text, then drop the resulting "quote prefix" formatting flag so the
# cells end up with no special style (matching the other quarter sheets).
$q1.Range("B2").Value = "'004223"
$q1.Range("C2").Value = "金信多策略精选灵活配置混合"
$q1.Range("D2").Value = "'0.36"
$q1.Range("E2").Value = "'93.14"
$q1.Range("F2").Value = "'4.79"
$q1.Range("G2").Value = "'0.0172"
$q1.Range("B2").ClearFormats()
$q1.Range("D2:G2").ClearFormats()
$q1.Range("H2").Value = 6

# --- Fill in the new "总计" sheet: headers + the 4 data rows (2022-Q1 on
# top, followed by the previously-existing quarters shifted down by one).
$totalWs = $newTotal
$totalWs.Range("B1").Value = "日期"
$totalWs.Range("C1").Value = "持有数量(只)"
$totalWs.Range("D1").Value = "持有市值(亿元)"

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 1
$totalWs.Range("D2").Value = 0.02

$totalWs.Range("A3").Value = 1
$totalWs.Range("B3").Value = "2021-Q4"
$totalWs.Range("C3").Value = 1
$totalWs.Range("D3").Value = 0.1

$totalWs.Range("A4").Value = 2
$totalWs.Range("B4").Value = "2021-Q3"
$totalWs.Range("C4").Value = 3
$totalWs.Range("D4").Value = 1.02

$totalWs.Range("A5").Value = 3
$totalWs.Range("B5").Value = "2021-Q2"
$totalWs.Range("C5").Value = 1
$totalWs.Range("D5").Value = 0

# Restore the original active sheet/tab (our edits above left "总计"
# selected; the workbook originally had "2021-Q2" as the active sheet).
$wb.Worksheets.Item("2021-Q2").Activate()
